$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row: add Wins / Losses / Ties in AD1:AF1, matching the style of
# the existing header cells (e.g. A1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AD1:AF1").Style = $ws.Range("A1").Style

# Data rows 2-46: team record values (same for every row).
$lastRow = 46
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 77   # AD
    $ws.Cells.Item($r, 31).Value = 85   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
